# Updated cryptos list values per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns being touched so Excel
# does not reinterpret text like "591.70" or "5.10" as numbers
# (which would silently drop the significant trailing zero).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.771.73"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.494.99"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "591.70"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "173.73"
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D9").Value = "2.493.93"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").Value = "5.10"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").Value = "26.29"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "2.951.38"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "67.613.06"
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").Value = "2.492.32"
$ws.Range("E18").Value = "  -2.20%  "
$ws.Range("D19").Value = "11.69"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("D21").Value = "365.20"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "71.40"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -6.61%  "
$ws.Range("D27").Value = "9.89"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "2.608.44"
$ws.Range("D30").Value = "0.0₃0965"
$ws.Range("D31").Value = "534.29"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").Value = "8.25"
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("E33").Value = "  -4.49%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -4.64%  "
$ws.Range("D37").Value = "158.80"
$ws.Range("E37").Value = "  +0.91%  "
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("D39").Value = "18.61"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "5.12"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "0.349"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  -1.59%  "
$ws.Range("D46").Value = "144.76"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "3.68"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "0.548"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").Value = "0.0₆0274"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").Value = "1.69"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("E51").Value = "  -1.89%  "
